$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for "Feria Lagunitas de Puerto
# Montt - Albahaca" that sorts in as the new row 98 (by date), pushing the
# former rows 98-143 down to 99-144. Insert a row at 98 to reproduce that
# shift, then populate the new row.
$ws.Rows("98").Insert()

# The market/region/category/quality descriptive columns (A,B,C,E,F,G,H,I,R)
# are identical across every record in this sheet, so carry them over from
# the row that now sits just below the new one.
$ws.Cells.Item(98, 1).Value2 = $ws.Cells.Item(99, 1).Value2
$ws.Cells.Item(98, 2).Value2 = $ws.Cells.Item(99, 2).Value2
$ws.Cells.Item(98, 3).Value2 = $ws.Cells.Item(99, 3).Value2
$ws.Cells.Item(98, 4).NumberFormat = $ws.Cells.Item(99, 4).NumberFormat
$ws.Cells.Item(98, 4).Value2 = 44845
$ws.Cells.Item(98, 5).Value2 = $ws.Cells.Item(99, 5).Value2
$ws.Cells.Item(98, 6).Value2 = $ws.Cells.Item(99, 6).Value2
$ws.Cells.Item(98, 7).Value2 = $ws.Cells.Item(99, 7).Value2
$ws.Cells.Item(98, 8).Value2 = $ws.Cells.Item(99, 8).Value2
$ws.Cells.Item(98, 9).Value2 = $ws.Cells.Item(99, 9).Value2
$ws.Cells.Item(98, 10).Value2 = 90
$ws.Cells.Item(98, 11).Value2 = 5000
$ws.Cells.Item(98, 12).Value2 = 5000
$ws.Cells.Item(98, 13).Value2 = 5000
$ws.Cells.Item(98, 14).Value2 = "$/paquete"
$ws.Cells.Item(98, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(98, 16).Value2 = 5000
$ws.Cells.Item(98, 17).Value2 = 1
$ws.Cells.Item(98, 18).Value2 = $ws.Cells.Item(99, 18).Value2
